$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.517.20"
$ws.Range("E2").Value = "  -1.27%  "

$ws.Range("D3").Value = "1.854.53"
$ws.Range("E3").Value = "  -0.36%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").Value = "'242.24"
$ws.Range("E5").Value = "  -0.91%  "

$ws.Range("D6").Value = "'0.6315"
$ws.Range("E6").Value = "  -4.18%  "

$ws.Range("D7").Value = "'1.002"
$ws.Range("E7").Value = "  +0.19%  "

$ws.Range("D8").Value = "'0.07580"
$ws.Range("E8").Value = "  -0.12%  "

$ws.Range("D9").Value = "'0.2989"
$ws.Range("E9").Value = "  -0.28%  "

$ws.Range("D10").Value = "'24.58"
$ws.Range("E10").Value = "  -0.61%  "

$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "'0.07714"
$ws.Range("E11").Value = "  +0.93%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.927.07"
$ws.Range("E12").Value = "  +3.31%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.018"
$ws.Range("E13").Value = "  -1.25%  "

$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.6914"
$ws.Range("E14").Value = "  -0.17%  "

$ws.Range("D15").Value = "'83.59"

$ws.Range("D16").Value = "'0.000009888"
$ws.Range("E16").Value = "  +1.66%  "

$ws.Range("D17").Value = "2.176.59"
$ws.Range("E17").Value = "  +2.77%  "

$ws.Range("D18").Value = "'6.217"
$ws.Range("E18").Value = "  +1.15%  "

$ws.Range("D19").Value = "29.676.56"
$ws.Range("E19").Value = "  -0.79%  "

$ws.Range("D20").Value = "'233.90"
$ws.Range("E20").Value = "  -1.38%  "

$ws.Range("D21").Value = "'12.59"
$ws.Range("E21").Value = "  -1.07%  "

$ws.Range("E22").Value = "  +0.12%  "

$ws.Range("D23").Value = "'7.691"
$ws.Range("E23").Value = "  -1.79%  "

$ws.Range("D24").Value = "'1.002"
$ws.Range("E24").Value = "  +0.06%  "

$ws.Range("D25").Value = "'155.60"
$ws.Range("E25").Value = "  -1.78%  "

$ws.Range("D26").Value = "'0.1398"
$ws.Range("E26").Value = "  -3.30%  "

$ws.Range("D27").Value = "'8.496"
$ws.Range("E27").Value = "  -1.44%  "

$ws.Range("D28").Value = "'17.74"
$ws.Range("E28").Value = "  -0.99%  "

$ws.Range("D29").Value = "'1.477"
$ws.Range("E29").Value = "  -1.28%  "

$ws.Range("D30").Value = "'0.05775"
$ws.Range("E30").Value = "  -4.59%  "

$ws.Range("D31").Value = "'1.257"
$ws.Range("E31").Value = "  -2.61%  "

$ws.Range("D32").Value = "'4.136"
$ws.Range("E32").Value = "  -0.83%  "

$ws.Range("D33").Value = "'4.032"
$ws.Range("E33").Value = "  -1.60%  "

$ws.Range("D34").Value = "'1.893"
$ws.Range("E34").Value = "  +0.87%  "

$ws.Range("D35").Value = "'1.171"
$ws.Range("E35").Value = "  -1.11%  "

$ws.Range("D36").Value = "'0.7215"
$ws.Range("E36").Value = "  -1.96%  "

$ws.Range("D37").Value = "'2.591"
$ws.Range("E37").Value = "  -0.71%  "

$ws.Range("D38").Value = "1.258.98"
$ws.Range("E38").Value = "  +4.09%  "

$ws.Range("D39").Value = "'2.807"
$ws.Range("E39").Value = "  -0.32%  "

$ws.Range("D40").Value = "'0.01808"
$ws.Range("E40").Value = "  +0.47%  "

$ws.Range("D41").Value = "'0.9041"
$ws.Range("E41").Value = "  -1.50%  "

$ws.Range("D42").Value = "'6.149"
$ws.Range("E42").Value = "  -2.63%  "

$ws.Range("D43").Value = "2.091.67"
$ws.Range("E43").Value = "  +3.13%  "

$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("D45").Value = "'67.91"
$ws.Range("E45").Value = "  +0.52%  "

$ws.Range("D46").Value = "'101.52"
$ws.Range("E46").Value = "  -0.12%  "

$ws.Range("D47").Value = "'7.338"
$ws.Range("E47").Value = "  -2.40%  "

$ws.Range("D48").Value = "'0.00000000119"
$ws.Range("E48").Value = "  -3.49%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.203"
$ws.Range("E49").Value = "  -0.03%  "

$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").Value = "'0.4046"
$ws.Range("E50").Value = "  -0.86%  "

$ws.Range("D51").Value = "'1.712"
$ws.Range("E51").Value = "  +1.35%  "
